# Adding one new search testcase to the "Test Cases" sheet.
#
# Summary of the edit:
#  - Almost every existing test case (rows 2-105) is switched from
#    Runmode "Y" to Runmode "N" (i.e. no longer executed), and picks up a
#    "SKIP" result where it did not already have one.
#  - The old combined "POSTs / ALL search results page" post record-view
#    test case (row 106, Jira OPQA-555) is reworked: its description is
#    expanded with a full list of field-level verification steps, its
#    Jira reference becomes "OPQA-555|OPQA-556", and it keeps
#    Runmode "Y" / Results "PASS".
#  - Column B is widened (best-fit) to accommodate the longer Jira id
#    text, and the view is scrolled/selected near the bottom of the
#    sheet where the edit happened.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Bulk re-mark existing test cases as not-run -------------------------
# D2:D105 -> "N"
$ws.Range("D2:D105").Value = "N"

# Rows 100-105 did not already carry a "Results" value (or had it as
# "PASS") - they now all read "SKIP".
$ws.Range("E100:E105").Value = "SKIP"

# --- Row 106: the reworked / expanded POSTs test case ---------------------
$nl = [char]10
$newDescription = "Verify that record view page of a post gets displayed when user clicks on article title in POSTs search results page" + $nl + `
  "Verify that following fields get displayed correctly for a post in record view page:" + $nl + `
  "a)Title" + $nl + `
  " b)Creation date and time " + $nl + `
  "c)Last edited date and time " + $nl + `
  "d)Author " + $nl + `
  "e)Author details " + $nl + `
  "f)Post content " + $nl + `
  "g)Likes count " + $nl + `
  "h)Comments count " + $nl + `
  "i)Views count"

$ws.Range("B106").Value = "OPQA-555|OPQA-556"
$ws.Range("C106").Value = $newDescription
$ws.Range("C106").WrapText = $true
$ws.Rows.Item(106).RowHeight = 165

# Row 106 keeps Runmode = Y and now has a Results = PASS
$ws.Range("D106").Value = "Y"
$ws.Range("E106").Value = "PASS"

# --- Column sizing ----------------------------------------------------
# Column B now gets its own (wider) best-fit width, separate from column A.
$ws.Columns.Item(2).ColumnWidth = 19.5

# --- View state: scroll near the bottom and select the last edited cell --
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 97
$win.ScrollColumn = 1
$ws.Range("D106").Select()

Write-Output "done"
